# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Final table (rows 16-30, cols B:G) after the update.
# Worker rows now interleave JOSE DAVID ALMEIDA LEONES / RUBEN DARIO GOMEZ OROZCO
# period by period (2305..2311), and several "Salario Basico" (G) / "Valor Mora" (F)
# amounts were refreshed along with the new statement data.
$data = @(
    @(16, "CC", "73198932",   "ALEXANDER GUZMAN RESTREPO", "2206", 160000, 2000000),
    @(17, "CC", "73006956",   "JOSE DAVID ALMEIDA LEONES",  "2305",  46400, 3000000),
    @(18, "CC", "1143351433", "RUBEN DARIO GOMEZ OROZCO",   "2305",  80000, 2000000),
    @(19, "CC", "73006956",   "JOSE DAVID ALMEIDA LEONES",  "2306",  46400, 3000000),
    @(20, "CC", "1143351433", "RUBEN DARIO GOMEZ OROZCO",   "2306",  80000, 2000000),
    @(21, "CC", "73006956",   "JOSE DAVID ALMEIDA LEONES",  "2307",  46400, 3000000),
    @(22, "CC", "1143351433", "RUBEN DARIO GOMEZ OROZCO",   "2307",  80000, 2000000),
    @(23, "CC", "73006956",   "JOSE DAVID ALMEIDA LEONES",  "2308", 120000, 3000000),
    @(24, "CC", "1143351433", "RUBEN DARIO GOMEZ OROZCO",   "2308",  80000, 2000000),
    @(25, "CC", "73006956",   "JOSE DAVID ALMEIDA LEONES",  "2309", 120000, 3000000),
    @(26, "CC", "1143351433", "RUBEN DARIO GOMEZ OROZCO",   "2309",  80000, 2000000),
    @(27, "CC", "73006956",   "JOSE DAVID ALMEIDA LEONES",  "2310", 120000, 3000000),
    @(28, "CC", "1143351433", "RUBEN DARIO GOMEZ OROZCO",   "2310",  80000, 2000000),
    @(29, "CC", "73006956",   "JOSE DAVID ALMEIDA LEONES",  "2311", 112000, 3000000),
    @(30, "CC", "1143351433", "RUBEN DARIO GOMEZ OROZCO",   "2311",  74667, 2000000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}
